# Add a "TYPE" column to the Mapping sheet: header in D1, then "TV" for the
# TV-platform rows (2-30) and "RADIO" for the RADIO-platform rows (31-38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("D1").Value = "TYPE"

for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 4).Value = "TV"
}

for ($r = 31; $r -le 38; $r++) {
    $ws.Cells.Item($r, 4).Value = "RADIO"
}

# Bring the Mapping sheet to the front and move the selection, matching the
# view state captured in the saved workbook.
$ws.Activate() | Out-Null
$ws.Range("F35").Select() | Out-Null
